$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the diff: (cellRef, newValue, forceText)
$updates = @(
    @('D2', '58.399.31', $False),
    @('E2', '  -4.12%  ', $False),
    @('D3', '2.534.25', $False),
    @('E3', '  -3.68%  ', $False),
    @('E4', '  +0.03%  ', $False),
    @('D5', '507.26', $True),
    @('E5', '  -4.18%  ', $False),
    @('D6', '144.16', $True),
    @('E6', '  -7.21%  ', $False),
    @('D7', '0.999', $True),
    @('E7', '  +0.07%  ', $False),
    @('D8', '0.564', $True),
    @('E8', '  -4.19%  ', $False),
    @('D9', '2.537.35', $False),
    @('E9', '  -3.83%  ', $False),
    @('D10', '6.09', $True),
    @('E10', '  -8.40%  ', $False),
    @('E11', '  -7.06%  ', $False),
    @('E12', '  -5.68%  ', $False),
    @('E13', '  -0.62%  ', $False),
    @('D14', '2.978.81', $False),
    @('E14', '  -3.73%  ', $False),
    @('D15', '58.394.22', $False),
    @('E15', '  -4.14%  ', $False),
    @('D16', '20.66', $True),
    @('E16', '  -6.18%  ', $False),
    @('E17', '  -6.39%  ', $False),
    @('D18', '2.534.28', $False),
    @('E18', '  -3.75%  ', $False),
    @('D19', '4.53', $True),
    @('E19', '  -5.06%  ', $False),
    @('D20', '334.55', $True),
    @('E20', '  -5.28%  ', $False),
    @('D21', '10.09', $True),
    @('E21', '  -4.98%  ', $False),
    @('D22', '0.998', $True),
    @('E22', '  -0.25%  ', $False),
    @('D23', '5.95', $True),
    @('E23', '  -4.89%  ', $False),
    @('D24', '60.57', $True),
    @('E24', '  -1.75%  ', $False),
    @('D25', '0.408', $True),
    @('E25', '  -4.97%  ', $False),
    @('D26', '0.999', $True),
    @('E26', '  -0.12%  ', $False),
    @('D27', '0.159', $True),
    @('E27', '  -5.62%  ', $False),
    @('D28', '2.648.68', $False),
    @('E28', '  -3.83%  ', $False),
    @('D29', '0.0₃0784', $False),
    @('E29', '  -9.48%  ', $False),
    @('D30', '6.95', $True),
    @('E30', '  -6.26%  ', $False),
    @('B32', 'Monero', $False),
    @('C32', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', $False),
    @('D32', '149.42', $True),
    @('E32', '  -1.04%  ', $False),
    @('B33', 'EthereumClassic', $False),
    @('C33', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', $False),
    @('D33', '18.53', $True),
    @('E33', '  -5.00%  ', $False),
    @('B34', 'Aptos', $False),
    @('C34', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', $False),
    @('D34', '5.82', $True),
    @('E34', '  -5.39%  ', $False),
    @('D35', '1.53', $True),
    @('E35', '  -5.52%  ', $False),
    @('D36', '0.916', $True),
    @('E36', '  +3.70%  ', $False),
    @('D37', '3.90', $True),
    @('E37', '  -6.05%  ', $False),
    @('E38', '  -7.53%  ', $False),
    @('D39', '35.95', $True),
    @('E39', '  -1.83%  ', $False),
    @('D40', '0.823', $True),
    @('E40', '  -11.65%  ', $False),
    @('D41', '1.39', $True),
    @('E41', '  -7.18%  ', $False),
    @('D42', '283.78', $True),
    @('E42', '  -7.20%  ', $False),
    @('E43', '  -7.91%  ', $False),
    @('D44', '0.0996', $True),
    @('E44', '  -2.69%  ', $False),
    @('D45', '0.997', $True),
    @('D46', '0.600', $True),
    @('E46', '  -6.47%  ', $False),
    @('D47', '0.0533', $True),
    @('E47', '  -5.28%  ', $False),
    @('B48', 'WhiteBITCoin', $False),
    @('C48', 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt', $False),
    @('D48', '10.29', $True),
    @('E48', '  -0.55%  ', $False),
    @('B49', 'EnergySwap', $False),
    @('C49', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', $False),
    @('D49', '18.59', $True),
    @('E49', '  -5.47%  ', $False),
    @('E50', '  -5.04%  ', $False),
    @('D51', '4.52', $True),
    @('E51', '  -8.57%  ', $False),
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $forceText = $u[2]
    $cell = $ws.Range($ref)
    if ($forceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

Write-Host "Applied $($updates.Count) cell updates."
